$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as plain text in the source workbook (even
# when they look numeric). Force those specific cells to keep a Text number
# format first, otherwise Excel would silently coerce the assigned string into
# a number and drop meaningful trailing zeros (e.g. "1.000" -> 1).
$textPriceCells = @('D4','D5','D6','D7','D8','D9','D10','D11','D13','D14','D16','D17','D19','D20','D21','D22','D23','D25','D26','D27','D28','D29','D30','D31','D32','D33','D34','D35','D36','D38','D39','D40','D41','D42','D43','D45','D46','D47','D48','D49','D50','D51')
foreach ($c in $textPriceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.037.83'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.833.89'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '241.72'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = '0.6275'
$ws.Range('E6').Value = '  -5.60%  '
$ws.Range('D7').Value = '0.9995'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.07665'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('D9').Value = '0.2919'
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').Value = '22.76'
$ws.Range('E10').Value = '  -2.92%  '
$ws.Range('D11').Value = '0.07736'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '1.830.72'
$ws.Range('E12').Value = '  -1.83%  '
$ws.Range('D13').Value = '4.954'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('D14').Value = '0.6660'
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('E15').Value = '  +18.37%  '
$ws.Range('D16').Value = '82.79'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '6.008'
$ws.Range('E17').Value = '  -3.11%  '
$ws.Range('D18').Value = '29.021.29'
$ws.Range('D19').Value = '225.72'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('D20').Value = '12.33'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('D21').Value = '0.9991'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '7.208'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').Value = '8.437'
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('D26').Value = '0.1370'
$ws.Range('E26').Value = '  -2.92%  '
$ws.Range('D27').Value = '17.89'
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('D28').Value = '1.494'
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('D29').Value = '4.067'
$ws.Range('E29').Value = '  -1.84%  '
$ws.Range('D30').Value = '4.030'
$ws.Range('E30').Value = '  -0.99%  '
$ws.Range('D31').Value = '1.197'
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('D32').Value = '0.05194'
$ws.Range('E32').Value = '  -2.71%  '
$ws.Range('D33').Value = '1.851'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('D34').Value = '0.7389'
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').Value = '1.145'
$ws.Range('E35').Value = '  -1.43%  '
$ws.Range('D36').Value = '2.695'
$ws.Range('D37').Value = '1.271.96'
$ws.Range('E37').Value = '  -3.90%  '
$ws.Range('D38').Value = '2.760'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = '0.01786'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').Value = '6.327'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('D41').Value = '0.8955'
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Value = '101.47'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('D44').Value = '1.979.61'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').Value = '0.00000000125'
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('D46').Value = '64.37'
$ws.Range('E46').Value = '  -2.14%  '
$ws.Range('D47').Value = '0.5118'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').Value = '0.4005'
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('D49').Value = '8.835'
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('D50').Value = '0.05763'
$ws.Range('D51').Value = '1.641'
$ws.Range('E51').Value = '  -6.76%  '
